# Lab4_RunningTimeStats.xlsx edit
# - widen columns A, D, H (prep for the new draft R-file content)
# - select the whole sheet (draft state left by the author before saving)
# - drop the AutoFilter that was on A1:M1 (sort order / filter no longer wanted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A, D, H to fit the new, longer content.
$ws.Columns.Item(1).ColumnWidth = 83.83333333333333
$ws.Columns.Item(4).ColumnWidth = 20.666666666666668
$ws.Columns.Item(8).ColumnWidth = 36

# Select the entire sheet (equivalent to Ctrl+A / clicking the corner selector).
[void]$ws.Cells.Select()

# Remove the AutoFilter dropdowns from the header row.
$ws.AutoFilterMode = $false
